$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10.82091761467928
$ws.Range("C2").Value = 8.544792735556168
$ws.Range("E2").Value = 13.68982799324578
$ws.Range("F2").Value = 49.17445699580077
$ws.Range("G2").Value = 3.724083580899947
$ws.Range("I2").Value = 32.34390339147611
$ws.Range("J2").Value = 10.04783658443079
$ws.Range("K2").Value = 12.19871911858792
$ws.Range("M2").Value = 16.87858594162278
$ws.Range("B3").Value = 10.69221808657815
$ws.Range("C3").Value = 8.402181951277045
$ws.Range("E3").Value = 13.72135786726382
$ws.Range("F3").Value = 48.89058132267645
$ws.Range("G3").Value = 3.727049664549777
$ws.Range("I3").Value = 32.20303204209402
$ws.Range("J3").Value = 10.05220735619478
$ws.Range("K3").Value = 12.11671602305847
$ws.Range("M3").Value = 16.87685127683544
$ws.Range("B4").Value = 10.61624489910634
$ws.Range("C4").Value = 8.315971919726122
$ws.Range("E4").Value = 13.74467724691989
$ws.Range("F4").Value = 48.72502639962015
$ws.Range("G4").Value = 3.728964763229378
$ws.Range("I4").Value = 32.11990677658027
$ws.Range("J4").Value = 10.05537503852897
$ws.Range("K4").Value = 12.06969799872791
$ws.Range("M4").Value = 16.8802433366178
$ws.Range("B5").Value = 10.58609354199783
$ws.Range("C5").Value = 8.281228882967403
$ws.Range("E5").Value = 13.75517485861391
$ws.Range("F5").Value = 48.65980473168558
$ws.Range("G5").Value = 3.729768881173714
$ws.Range("I5").Value = 32.08689308786929
$ws.Range("J5").Value = 10.05678757681297
$ws.Range("K5").Value = 12.05139447732544
$ws.Range("M5").Value = 16.88274799415289
$ws.Range("B6").Value = 10.58113693564802
$ws.Range("C6").Value = 8.275484666574179
$ws.Range("E6").Value = 13.75697802620585
$ws.Range("F6").Value = 48.64911143221967
$ws.Range("G6").Value = 3.729903838243018
$ws.Range("I6").Value = 32.08146352333415
$ws.Range("J6").Value = 10.05702947613073
$ws.Range("K6").Value = 12.04840745304734
$ws.Range("M6").Value = 16.88323171091791
$ws.Range("B7").Value = 10.61583494075203
$ws.Range("C7").Value = 8.315501728090515
$ws.Range("E7").Value = 13.74481479511009
$ws.Range("F7").Value = 48.72413765888356
$ws.Range("G7").Value = 3.728975511775121
$ws.Range("I7").Value = 32.11945804186879
$ws.Range("J7").Value = 10.0553935958697
$ws.Range("K7").Value = 12.06944765862416
$ws.Range("M7").Value = 16.88027256935865
$ws.Range("B8").Value = 10.77593476890722
$ws.Range("C8").Value = 8.495368458445023
$ws.Range("E8").Value = 13.69987712156025
$ws.Range("F8").Value = 49.07478894527618
$ws.Range("G8").Value = 3.725086846470136
$ws.Range("I8").Value = 32.29463957483985
$ws.Range("J8").Value = 10.04924315496381
$ws.Range("K8").Value = 12.16976459567461
$ws.Range("M8").Value = 16.87706390395017
$ws.Range("B9").Value = 11.1121698121545
$ws.Range("C9").Value = 8.856710905672585
$ws.Range("E9").Value = 13.6432053219045
$ws.Range("F9").Value = 49.82971865173432
$ws.Range("G9").Value = 3.718202451527256
$ws.Range("I9").Value = 32.66429779204037
$ws.Range("J9").Value = 10.04102405113513
$ws.Range("K9").Value = 12.39202732242138
$ws.Range("M9").Value = 16.90603795925323
$ws.Range("B10").Value = 11.37014178046914
$ws.Range("C10").Value = 9.124632379941907
$ws.Range("E10").Value = 13.62076440326372
$ws.Range("F10").Value = 50.42238885041401
$ws.Range("G10").Value = 3.713590947596527
$ws.Range("I10").Value = 32.95093435079709
$ws.Range("J10").Value = 10.03733006999926
$ws.Range("K10").Value = 12.5695979132339
$ws.Range("M10").Value = 16.94864731662132
$ws.Range("B11").Value = 11.48930204034355
$ws.Range("C11").Value = 9.246466898119028
$ws.Range("E11").Value = 13.61472105614144
$ws.Range("F11").Value = 50.69956756056555
$ws.Range("G11").Value = 3.71158884010067
$ws.Range("I11").Value = 33.08439333640776
$ws.Range("J11").Value = 10.03615914895285
$ws.Range("K11").Value = 12.65317907093557
$ws.Range("M11").Value = 16.97260818085384
$ws.Range("B12").Value = 11.5346369673955
$ws.Range("C12").Value = 9.292549517229229
$ws.Range("E12").Value = 13.61303065647377
$ws.Range("F12").Value = 50.80555381324127
$ws.Range("G12").Value = 3.710844365425079
$ws.Range("I12").Value = 33.13535261134326
$ws.Range("J12").Value = 10.03578903026289
$ws.Range("K12").Value = 12.68520571325775
$ws.Range("M12").Value = 16.98233406572315
$ws.Range("B13").Value = 11.5248646367723
$ws.Range("C13").Value = 9.282627909495208
$ws.Range("E13").Value = 13.61336813221394
$ws.Range("F13").Value = 50.78268324172787
$ws.Range("G13").Value = 3.711004094214453
$ws.Range("I13").Value = 33.12435920678642
$ws.Range("J13").Value = 10.03586548239385
$ws.Range("K13").Value = 12.67829190249715
$ws.Range("M13").Value = 16.98021050641026
$ws.Range("B14").Value = 11.4930278080943
$ws.Range("C14").Value = 9.250259477924896
$ws.Range("E14").Value = 13.6145700057654
$ws.Range("F14").Value = 50.70826694986989
$ws.Range("G14").Value = 3.71152731800996
$ws.Range("I14").Value = 33.08857745987579
$ws.Range("J14").Value = 10.03612723026308
$ws.Range("K14").Value = 12.65580653717359
$ws.Range("M14").Value = 16.97339529574281
$ws.Range("B15").Value = 11.47355295059069
$ws.Range("C15").Value = 9.230424548355137
$ws.Range("E15").Value = 13.61538404404604
$ws.Range("F15").Value = 50.66281633792882
$ws.Range("G15").Value = 3.711849586865047
$ws.Range("I15").Value = 33.0667143674016
$ws.Range("J15").Value = 10.03629710255645
$ws.Range("K15").Value = 12.64208179504756
$ws.Range("M15").Value = 16.96930556160498
$ws.Range("B16").Value = 11.36238660093681
$ws.Range("C16").Value = 9.116665250210138
$ws.Range("E16").Value = 13.62124309928975
$ws.Range("F16").Value = 50.40442148481563
$ws.Range("G16").Value = 3.71372370868858
$ws.Range("I16").Value = 32.94227253286429
$ws.Range("J16").Value = 10.03741684652209
$ws.Range("K16").Value = 12.56418981981415
$ws.Range("M16").Value = 16.94717305010506
$ws.Range("B17").Value = 11.29461801862241
$ws.Range("C17").Value = 9.046832203040701
$ws.Range("E17").Value = 13.62590374840529
$ws.Range("F17").Value = 50.24779869404065
$ws.Range("G17").Value = 3.71489787289171
$ws.Range("I17").Value = 32.86670315647255
$ws.Range("J17").Value = 10.03823428192665
$ws.Range("K17").Value = 12.5171047533687
$ws.Range("M17").Value = 16.93476400364911
$ws.Range("B18").Value = 11.25581233186818
$ws.Range("C18").Value = 9.006664627151146
$ws.Range("E18").Value = 13.62897654963247
$ws.Range("F18").Value = 50.15842989998591
$ws.Range("G18").Value = 3.715582232520501
$ws.Range("I18").Value = 32.82352752502359
$ws.Range("J18").Value = 10.03875240399475
$ws.Range("K18").Value = 12.49028819528267
$ws.Range("M18").Value = 16.92805791591607
$ws.Range("B19").Value = 11.24270458988032
$ws.Range("C19").Value = 8.993065789791387
$ws.Range("E19").Value = 13.63008431372636
$ws.Range("F19").Value = 50.12829613286191
$ws.Range("G19").Value = 3.715815495096443
$ws.Range("I19").Value = 32.80895937161185
$ws.Range("J19").Value = 10.03893606679384
$ws.Range("K19").Value = 12.48125498072354
$ws.Range("M19").Value = 16.92586158597608
$ws.Range("B20").Value = 11.30181454632226
$ws.Range("C20").Value = 9.054266547396688
$ws.Range("E20").Value = 13.62536703650931
$ws.Range("F20").Value = 50.26439775982126
$ws.Range("G20").Value = 3.714771949041641
$ws.Range("I20").Value = 32.87471778501679
$ws.Range("J20").Value = 10.03814230123102
$ws.Range("K20").Value = 12.52208975049169
$ws.Range("M20").Value = 16.93604037325651
$ws.Range("B21").Value = 11.5023736987054
$ws.Range("C21").Value = 9.259768688491212
$ws.Range("E21").Value = 13.61420076333494
$ws.Range("F21").Value = 50.73009752760876
$ws.Range("G21").Value = 3.711373263862412
$ws.Range("I21").Value = 33.09907615471889
$ws.Range("J21").Value = 10.03604835954844
$ws.Range("K21").Value = 12.6624010393739
$ws.Range("M21").Value = 16.97537943186562
$ws.Range("B22").Value = 11.63466162544358
$ws.Range("C22").Value = 9.393744026638636
$ws.Range("E22").Value = 13.61038837289211
$ws.Range("F22").Value = 51.04040098895357
$ws.Range("G22").Value = 3.709231727227666
$ws.Range("I22").Value = 33.24815300958318
$ws.Range("J22").Value = 10.0351070052449
$ws.Range("K22").Value = 12.75628289688825
$ws.Range("M22").Value = 17.00488938813142
$ws.Range("B23").Value = 11.56396210001082
$ws.Range("C23").Value = 9.322284348791705
$ws.Range("E23").Value = 13.61210458848404
$ws.Range("F23").Value = 50.87426411427939
$ws.Range("G23").Value = 3.710367439330235
$ws.Range("I23").Value = 33.16837082366231
$ws.Range("J23").Value = 10.03557033411208
$ws.Range("K23").Value = 12.70598591595296
$ws.Range("M23").Value = 16.98879382396921
$ws.Range("B24").Value = 11.29856050848684
$ws.Range("C24").Value = 9.050905537847335
$ws.Range("E24").Value = 13.62560845873012
$ws.Range("F24").Value = 50.2568912124571
$ws.Range("G24").Value = 3.714828850186652
$ws.Range("I24").Value = 32.87109352823872
$ws.Range("J24").Value = 10.03818373567002
$ws.Range("K24").Value = 12.51983524233976
$ws.Range("M24").Value = 16.93546199283935
$ws.Range("B25").Value = 11.01909398834256
$ws.Range("C25").Value = 8.758328074466375
$ws.Range("E25").Value = 13.65516574988953
$ws.Range("F25").Value = 49.61858348625354
$ws.Range("G25").Value = 3.719986068541739
$ws.Range("I25").Value = 32.34390339147611
$ws.Range("J25").Value = 10.04283595407114
$ws.Range("K25").Value = 12.32929403495316
$ws.Range("M25").Value = 16.8944395330318
